$wb = $excel.ActiveWorkbook

# --- Sheet 1: Means ---
$ws1 = $wb.Worksheets.Item("Means")

# Header row: add new columns F and G
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# Row 2 (% White)
$ws1.Range("F2").Value = 52
$ws1.Range("G2").Value = 50

# Row 3 (% Black or African American)
$ws1.Range("F3").Value = 10
$ws1.Range("G3").Value = 9.4

# Row 4 (% Other)
$ws1.Range("F4").Value = 38
$ws1.Range("G4").Value = 40

# Row 5 (% Hispanic)
$ws1.Range("F5").Value = 50
$ws1.Range("G5").Value = 45

# Row 6 (Median Income)
$ws1.Range("F6").Value = 55
$ws1.Range("G6").Value = 62

# Row 7 (% Below Poverty Line)
$ws1.Range("F7").Value = 11
$ws1.Range("G7").Value = 9.8

# Row 8 (% Below Half the Poverty Line)
$ws1.Range("F8").Value = 8
$ws1.Range("G8").Value = 7

# Row 9 (Total Cancer Risk) - also updates existing B-E values
$ws1.Range("B9").Value = 29
$ws1.Range("C9").Value = 31
$ws1.Range("D9").Value = 30
$ws1.Range("E9").Value = 30
$ws1.Range("F9").Value = 30
$ws1.Range("G9").Value = 30

# Row 10 (Total Respiratory) - also updates existing B-E values
$ws1.Range("B10").Value = 0.37
$ws1.Range("C10").Value = 0.43
$ws1.Range("D10").Value = 0.49
$ws1.Range("E10").Value = 0.5
$ws1.Range("F10").Value = 0.47
$ws1.Range("G10").Value = 0.45

# --- Sheet 2: Standard Deviations ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# Header row: add new columns F and G
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# Row 2 (% White)
$ws2.Range("F2").Value = 18
$ws2.Range("G2").Value = 19

# Row 3 (% Black or African American)
$ws2.Range("F3").Value = 10
$ws2.Range("G3").Value = 9.5

# Row 4 (% Other)
$ws2.Range("F4").Value = 15
$ws2.Range("G4").Value = 16

# Row 5 (% Hispanic)
$ws2.Range("F5").Value = 21
$ws2.Range("G5").Value = 21

# Row 6 (Median Income)
$ws2.Range("F6").Value = 22
$ws2.Range("G6").Value = 27

# Row 7 (% Below Poverty Line)
$ws2.Range("F7").Value = 11
$ws2.Range("G7").Value = 10

# Row 8 (% Below Half the Poverty Line)
$ws2.Range("F8").Value = 7.3
$ws2.Range("G8").Value = 6.9

# Row 9 (Total Cancer Risk) - also updates existing B-E values
$ws2.Range("B9").Value = 10
$ws2.Range("C9").Value = 8.9
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 0

# Row 10 (Total Respiratory) - also updates existing B-E values
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.23
$ws2.Range("D10").Value = 0.036
$ws2.Range("E10").Value = 0.015
$ws2.Range("F10").Value = 0.041
$ws2.Range("G10").Value = 0.05
